$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:C to B:D by inserting a new column at A
$ws.Columns.Item(1).Insert()

# New header for the inserted column
$ws.Range("B1").Value = "segments"
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122) # xlPasteFormats

# Copy the style of the (now shifted) name column onto the new index column
$ws.Range("B2:B20").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122) # xlPasteFormats

# Fill the new index column (0-based row index) with numeric values
for ($i = 2; $i -le 20; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 2
}

$excel.CutCopyMode = $false
